# Update "想去人数" (interest count) figures in the F column on the
# "展览" (Exhibition) and "全部类型" (All types) sheets to reflect the
# latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value  = 2827
$wsExhibit.Range("F10").Value = 1587
$wsExhibit.Range("F20").Value = 221
$wsExhibit.Range("F22").Value = 4
$wsExhibit.Range("F24").Value = 242
$wsExhibit.Range("F25").Value = 37
$wsExhibit.Range("F27").Value = 1802
$wsExhibit.Range("F29").Value = 431
$wsExhibit.Range("F30").Value = 98
$wsExhibit.Range("F34").Value = 462

# Sheet "全部类型" (sheet4.xml) - same events, rows shifted by +1
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value  = 2827
$wsAll.Range("F11").Value = 1587
$wsAll.Range("F21").Value = 221
$wsAll.Range("F23").Value = 4
$wsAll.Range("F25").Value = 242
$wsAll.Range("F26").Value = 37
$wsAll.Range("F28").Value = 1802
$wsAll.Range("F30").Value = 431
$wsAll.Range("F31").Value = 98
$wsAll.Range("F35").Value = 462
